$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update login password values (shared string contents change)
$ws.Range("B5").Value = "Fosroc@3"
$ws.Range("B11").Value = "Fosroc@0"
$ws.Range("B7").Value = "Fosroc@7"

# Update the selected cell / active cell on the sheet view (new xpath for login btn)
$ws.Range("F7").Select()
